$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A4 text (engineering thesis note expanded) and its row height
$ws.Range("A4").Value = "Inżynierka (opisac dane, potem model, znowu dane, i tabele GminaRelacja i jej triggery i czemu triggery)"
$ws.Rows.Item(4).RowHeight = 30.75

# Remove the now-obsolete TO_ASK block (rows 6-11)
$ws.Range("A6:A11").EntireRow.Delete()

# Update selection to mirror the saved cursor position
$ws.Range("A18").Select()
